$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Low-grade glioma")
$ws.Range("C2").Value = 0.083541683621774
$ws.Range("C3").Value = 0.0204091723342268
$ws.Range("C4").Value = 0.86626745181243
$ws.Range("C5").Value = 0.588174611595652
$ws.Range("C6").Value = 0.276767428792189
$ws.Range("C7").Value = 0.721643318278074
$ws.Range("C8").Value = 0.326402980186976
$ws.Range("C9").Value = 0.464613053587257

$ws = $wb.Worksheets.Item("Medulloblastoma")
$ws.Range("C2").Value = 0.717032215421925
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 0.987148764481946
$ws.Range("C5").Value = 0.971304078866354
$ws.Range("C6").Value = 0.638240668305458
$ws.Range("C7").Value = 0.0215000581483138
$ws.Range("C8").Value = 0.0136864060128698
$ws.Range("C9").Value = 0.824486303633589

$ws = $wb.Worksheets.Item("Mixed neuronal-glial tumor")
$ws.Range("C2").Value = 0.462739122413669
$ws.Range("C3").Value = 0.0578532713191215
$ws.Range("C4").Value = 0.999999999999986
$ws.Range("C5").Value = 0.656580631831596
$ws.Range("C6").Value = 0.999999999999963
$ws.Range("C7").Value = 0.336867791184883
$ws.Range("C8").Value = 0.499350846186132
$ws.Range("C9").Value = 0.943625773616395

$ws = $wb.Worksheets.Item("Ependymoma")
$ws.Range("C2").Value = 0.220196207956647
$ws.Range("C4").Value = 0.72429431453387
$ws.Range("C5").Value = 0.450754898244027
$ws.Range("C6").Value = 0.830379412511969
$ws.Range("C7").Value = 0.169620654830765
$ws.Range("C8").Value = 0.146688186321032
$ws.Range("C9").Value = 0.872589914881906

$ws = $wb.Worksheets.Item("Other high-grade glioma")
$ws.Range("C2").Value = 0.0549812456875758
$ws.Range("C3").Value = 0.0580535500081479
$ws.Range("C4").Value = 0.398228106546933
$ws.Range("C5").Value = 0.618892857690358
$ws.Range("C6").Value = 0.55188821452344
$ws.Range("C7").Value = 0.516918809678844
$ws.Range("C8").Value = 0.886473974807497
$ws.Range("C9").Value = 0.967326029074045

$ws = $wb.Worksheets.Item("Craniopharyngioma")
$ws.Range("C2").Value = 0.59919028340081
$ws.Range("C3").Value = 0.255060728744939
$ws.Range("C4").Value = 0.999999999999995
$ws.Range("C5").Value = 0.715504978662871
$ws.Range("C7").Value = 0.907258064516129
$ws.Range("C8").Value = 0.326612903225806
$ws.Range("C9").Value = 0.229161113594648

$ws = $wb.Worksheets.Item("ATRT")
$ws.Range("C2").Value = 0.432655369039808
$ws.Range("C3").Value = 0.694584286803966
$ws.Range("C4").Value = 0.455413419532409
$ws.Range("C5").Value = 0.169836296426685
$ws.Range("C6").Value = 0.511794871794871
$ws.Range("C7").Value = 0.219796130642813
$ws.Range("C8").Value = 0.28717351298587
$ws.Range("C9").Value = 0.288619115786948

$ws = $wb.Worksheets.Item("Meningioma")
$ws.Range("C9").Value = 0.094779946136451

$ws = $wb.Worksheets.Item("DIPG or DMG")
$ws.Range("C2").Value = 0.515458449661513
$ws.Range("C5").Value = 0.844026805515887
$ws.Range("C7").Value = 0.449023462236327
$ws.Range("C8").Value = 0.124542124542125
$ws.Range("C9").Value = 0.830698294656258

$ws = $wb.Worksheets.Item("Mesenchymal tumor")
$ws.Range("C9").Value = 0.924520737821253

$ws = $wb.Worksheets.Item("Neurofibroma plexiform")
$ws.Range("C9").Value = 0.646990031211512

$ws = $wb.Worksheets.Item("Non-neoplastic tumor")
$ws.Range("C9").Value = 0.125156485427265

$ws = $wb.Worksheets.Item("Germ cell tumor")
$ws.Range("C2").Value = 0.148251748251748
$ws.Range("C4").Value = 0.647552447552448
$ws.Range("C5").Value = 1
$ws.Range("C7").Value = 0.367676767676768
$ws.Range("C8").Value = 0.214141414141414
$ws.Range("C9").Value = 0.428523908344215

$ws = $wb.Worksheets.Item("Schwannoma")
$ws.Range("C9").Value = 0.0600624085865772

$ws = $wb.Worksheets.Item("Choroid plexus tumor")
$ws.Range("C2").Value = 0.645098039215686
$ws.Range("C3").Value = 0.119281045751634
$ws.Range("C6").Value = 0.999999999999997
$ws.Range("C7").Value = 0.903408211942422
$ws.Range("C8").Value = 0.261538461538462
$ws.Range("C9").Value = 0.283319589457907

$ws = $wb.Worksheets.Item("Other tumor")
$ws.Range("C2").Value = 0.100858857941891
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 0.490904927043123
$ws.Range("C5").Value = 0.75044100628083
$ws.Range("C6").Value = 0.665568860098732
$ws.Range("C7").Value = 0.930353912409497
$ws.Range("C8").Value = 0.838406383109114
$ws.Range("C9").Value = 0.596430178913545
